$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update task status values ("correccion anulacion de ots")
$ws.Range("B50").Value = "en proceso"
$ws.Range("B51").Value = "terminado"
$ws.Range("B53").Value = "terminado"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("A54").Select()
